# Insert a new data row at row 66 ("Magnum" variety, Región de Arica y
# Parinacota) which pushes the existing rows 66-171 down to 67-172, and
# fill in the new row's values (mirrors the xlsx diff: dimension grows
# from A1:R171 to A1:R172, one new record inserted mid-table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 66 (and everything below it) down by one row.
$ws.Range("A66").EntireRow.Insert()

# Populate the newly inserted row 66 with the new record's data.
$ws.Range("A66").Value = 5
$ws.Range("B66").Value = "Macroferia Regional de Talca"
$ws.Range("C66").Value = "Maule"
$ws.Range("D66").Value = 44791
$ws.Range("E66").Value = 7
$ws.Range("F66").Value = 100112031
$ws.Range("G66").Value = "Poroto verde"
$ws.Range("H66").Value = "Magnum"
$ws.Range("I66").Value = "Primera"
$ws.Range("J66").Value = 150
$ws.Range("K66").Value = 35000
$ws.Range("L66").Value = 35000
$ws.Range("M66").Value = 35000
$ws.Range("N66").Value = "`$/malla 25 kilos"
$ws.Range("O66").Value = "Región de Arica y Parinacota"
$ws.Range("P66").Value = 1400
$ws.Range("Q66").Value = 25
$ws.Range("R66").Value = "Hortaliza"
